$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# The "CreateDate" (row 15) and "LastUpdate" (row 17) fields were both
# documented with type "DATE"; change both to "TIMESTAMP".
$ws.Range("D15").Value = "TIMESTAMP"
$ws.Range("D17").Value = "TIMESTAMP"

# Reflect the final active selection used when the change was made.
$ws.Range("D17").Select()
